$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 312, shifting existing rows 312:378 down to 313:379.
$ws.Rows.Item(312).Insert()

# Populate the newly inserted row 312 with the new weekly record.
$ws.Range("A312").Value = 6
$ws.Range("B312").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C312").Value = 'Metropolitana'
$ws.Range("D312").Value = 45244
$ws.Range("E312").Value = 13
$ws.Range("F312").Value = 100112029
$ws.Range("G312").Value = 'Orégano'
$ws.Range("H312").Value = 'Sin especificar'
$ws.Range("I312").Value = 'Primera'
$ws.Range("J312").Value = 36
$ws.Range("K312").Value = 16000
$ws.Range("L312").Value = 16000
$ws.Range("M312").Value = 16000
$ws.Range("N312").Value = '$/docena de atados'
$ws.Range("O312").Value = 'Región Metropolitana'
$ws.Range("P312").Value = 5333
$ws.Range("Q312").Value = 3
$ws.Range("R312").Value = 'Hortaliza'
